# Refresh cryptos price/volume snapshot (and fix the ARBITRUM/HuobiToken row order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.543.20"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.915.78"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "'326.14"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "'0.4077"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "'0.08163"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'1.013"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'23.40"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").Value = "1.925.66"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'6.011"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'7.134"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'90.38"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'0.06801"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'17.71"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "29.559.11"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'5.625"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'11.80"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "'2.181"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "2.147.54"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'155.55"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "'6.487"
$ws.Range("E27").Value = "  +8.96%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "'2.102"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "'119.78"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "'1.032"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'5.523"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.397"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.559"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'0.02271"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "'0.06103"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'1.180"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'10.80"
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "'7.985"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "'0.1858"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").Value = "'2.491"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").Value = "'1.284"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'0.07718"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "'0.5573"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'1.949"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "'115.66"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "'72.84"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  +2.03%  "
